$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 84, shifting existing rows 84:169 down to 85:170.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record.
$ws.Range("A84").Value = 7
$ws.Range("B84").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C84").Value = "Ñuble"
$ws.Range("D84").Value = 44781
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 100112045
$ws.Range("G84").Value = "Zapallo"
$ws.Range("H84").Value = "Camote"
$ws.Range("I84").Value = "1a (guarda)"
$ws.Range("J84").Value = 240
$ws.Range("K84").Value = 700
$ws.Range("L84").Value = 800
$ws.Range("M84").Value = 750
$ws.Range("N84").Value = "$/kilo (volumen en unidades)"
$ws.Range("O84").Value = "Región de O'Higgins"
$ws.Range("P84").Value = 750
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"
